$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-12 hold 4 observation records. The update rotates the record
# contents between the rows (same underlying data, re-assigned to
# different Id/row slots): row 9 receives what was row 10's data, row 10
# receives what was row 11's data, row 11 receives what was row 12's
# data, and row 12 receives what was the original row 9's data. Columns
# C, I, P, S, T, U, V, W, Z, AB, AD, AE, AG, AT, AW, AX, AY are identical
# across these rows already, so only the columns below actually change.
$cols = @("A","B","D","E","F","G","H","Q","R","Y","AA")
# Columns that hold plain-text values which look like dates (must be
# re-entered with a leading apostrophe so Excel keeps them as text
# instead of silently converting them to date serials).
$textDateCols = @("Y","AA")
$rows = @(9,10,11,12)

# Snapshot all the current values first so the rotation reads are not
# affected by the writes that follow.
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Row 9 <- old Row 10, Row 10 <- old Row 11, Row 11 <- old Row 12, Row 12 <- old Row 9
$mapping = @{ 9 = 10; 10 = 11; 11 = 12; 12 = 9 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $val = $orig[$src][$c]
        if ($textDateCols -contains $c) {
            $ws.Range("$c$r").Value = "'" + $val
        } else {
            $ws.Range("$c$r").Value = $val
        }
    }
}
